# Auto-generated edit script applying the cryptos-list data refresh
# (GitHub Actions scheduled update) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.255.28"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "1.617.65"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'211.93"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'18.72"
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "1.843.67"
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("D13").Value = "1.623.18"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "26.273.96"
$ws.Range("D17").Value = "'62.20"
$ws.Range("E17").Value = "  +3.81%  "
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'200.65"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("E24").Value = "  +3.29%  "
$ws.Range("D25").Value = "'144.11"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'0.119"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  +8.95%  "
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.49"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.41"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").Value = "1.174.58"
$ws.Range("E36").Value = "  +4.84%  "
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'2.31"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'0.493"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").Value = "'0.791"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").Value = "'5.32"
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("D44").Value = "1.754.65"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "'92.53"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "0.0₆0105"
$ws.Range("E46").Value = "  +14.81%  "
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("D48").Value = "'53.56"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("E51").Value = "  -0.15%  "
